# Adds a new column F ("Ncompost(kg/kgTM)") computed as (0.008/0.6)-C for rows 2..93,
# and repositions/enlarges the existing chart to make room for the new column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column
$ws.Range("F1").Value = "Ncompost(kg/kgTM)"

# Row 2: explicit (non-shared) formula
$ws.Range("F2").Formula = "=(0.008/0.6)-C2"

# Rows 3..66 and 67..93: same formula pattern, filled down as two shared-formula
# groups (matching the existing break used by columns C/D/E at row 66/67).
$ws.Range("F3:F66").Formula = "=(0.008/0.6)-C3"
$ws.Range("F67:F93").Formula = "=(0.008/0.6)-C67"

# Move/resize the chart that sits on the sheet to make room for the new column
$chart = $ws.Shapes.Item(1)
$chart.TopLeftCell = $ws.Cells.Item(7, 10)
$chart.Left = $ws.Cells.Item(7, 10).Left + 492125
$chart.Top = $ws.Cells.Item(7, 10).Top + 9525
$bottomRightCell = $ws.Cells.Item(21, 16)
$chart.Width = ($bottomRightCell.Left + 492125) - $chart.Left
$chart.Height = ($bottomRightCell.Top + 171450) - $chart.Top

# Update the visible/selected cell state
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("D3").Select()
